$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Range("A2").Value = 'Última actualización: 20:11:44'
$ws1.Range("A3").Value = 'Total filas: 476'

# Row data updates (reshuffled schedule entries + newly scraped rows)
$ws1.Cells.Item(50, 1).Value = '07:12:53'
$ws1.Cells.Item(50, 2).Value = '08:01'
$ws1.Cells.Item(50, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(50, 4).Value = 49
$ws1.Cells.Item(50, 5).Value = 'LP1912'
$ws1.Cells.Item(51, 1).Value = '06:45:50'
$ws1.Cells.Item(51, 2).Value = '08:01'
$ws1.Cells.Item(51, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(51, 4).Value = 76
$ws1.Cells.Item(51, 5).Value = 'LP1912'
$ws1.Cells.Item(64, 1).Value = '08:29:19'
$ws1.Cells.Item(64, 2).Value = '08:29'
$ws1.Cells.Item(64, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(64, 4).Value = 0
$ws1.Cells.Item(64, 5).Value = 'LP1912'
$ws1.Cells.Item(65, 1).Value = '06:45:50'
$ws1.Cells.Item(65, 2).Value = '08:29'
$ws1.Cells.Item(65, 3).Value = '14_ABASTO'
$ws1.Cells.Item(65, 4).Value = 104
$ws1.Cells.Item(65, 5).Value = 'LP1912'
$ws1.Cells.Item(107, 1).Value = '09:21:49'
$ws1.Cells.Item(107, 2).Value = '10:03'
$ws1.Cells.Item(107, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(107, 4).Value = 42
$ws1.Cells.Item(107, 5).Value = 'LP1912'
$ws1.Cells.Item(108, 1).Value = '08:11:27'
$ws1.Cells.Item(108, 2).Value = '10:03'
$ws1.Cells.Item(108, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(108, 4).Value = 112
$ws1.Cells.Item(108, 5).Value = 'LP1912'
$ws1.Cells.Item(123, 1).Value = '10:04:17'
$ws1.Cells.Item(123, 2).Value = '10:25'
$ws1.Cells.Item(123, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(123, 4).Value = 21
$ws1.Cells.Item(123, 5).Value = 'LP1912'
$ws1.Cells.Item(124, 1).Value = '10:04:17'
$ws1.Cells.Item(124, 2).Value = '10:25'
$ws1.Cells.Item(124, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(124, 4).Value = 21
$ws1.Cells.Item(124, 5).Value = 'LP1912'
$ws1.Cells.Item(183, 1).Value = '12:11:45'
$ws1.Cells.Item(183, 2).Value = '12:11'
$ws1.Cells.Item(183, 3).Value = '15_ABASTO'
$ws1.Cells.Item(183, 4).Value = 0
$ws1.Cells.Item(183, 5).Value = 'LP1912'
$ws1.Cells.Item(184, 1).Value = '12:11:45'
$ws1.Cells.Item(184, 2).Value = '12:11'
$ws1.Cells.Item(184, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(184, 4).Value = 0
$ws1.Cells.Item(184, 5).Value = 'LP1912'
$ws1.Cells.Item(226, 1).Value = '11:47:13'
$ws1.Cells.Item(226, 2).Value = '13:11'
$ws1.Cells.Item(226, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(226, 4).Value = 84
$ws1.Cells.Item(226, 5).Value = 'LP1912'
$ws1.Cells.Item(227, 1).Value = '11:34:25'
$ws1.Cells.Item(227, 2).Value = '13:11'
$ws1.Cells.Item(227, 3).Value = '215_ALUAR'
$ws1.Cells.Item(227, 4).Value = 97
$ws1.Cells.Item(227, 5).Value = 'LP1912'
$ws1.Cells.Item(271, 1).Value = '13:51:48'
$ws1.Cells.Item(271, 2).Value = '14:25'
$ws1.Cells.Item(271, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(271, 4).Value = 34
$ws1.Cells.Item(271, 5).Value = 'LP1912'
$ws1.Cells.Item(272, 1).Value = '12:53:14'
$ws1.Cells.Item(272, 2).Value = '14:25'
$ws1.Cells.Item(272, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(272, 4).Value = 92
$ws1.Cells.Item(272, 5).Value = 'LP1912'
$ws1.Cells.Item(301, 1).Value = '13:39:24'
$ws1.Cells.Item(301, 2).Value = '15:25'
$ws1.Cells.Item(301, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(301, 4).Value = 106
$ws1.Cells.Item(301, 5).Value = 'LP1912'
$ws1.Cells.Item(302, 1).Value = '13:39:24'
$ws1.Cells.Item(302, 2).Value = '15:25'
$ws1.Cells.Item(302, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(302, 4).Value = 106
$ws1.Cells.Item(302, 5).Value = 'LP1912'
$ws1.Cells.Item(376, 1).Value = '17:34:55'
$ws1.Cells.Item(376, 2).Value = '18:03'
$ws1.Cells.Item(376, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(376, 4).Value = 29
$ws1.Cells.Item(376, 5).Value = 'LP1912'
$ws1.Cells.Item(377, 1).Value = '16:13:19'
$ws1.Cells.Item(377, 2).Value = '18:03'
$ws1.Cells.Item(377, 3).Value = '17_ROMERO'
$ws1.Cells.Item(377, 4).Value = 110
$ws1.Cells.Item(377, 5).Value = 'LP1912'
$ws1.Cells.Item(392, 1).Value = '17:34:55'
$ws1.Cells.Item(392, 2).Value = '18:33'
$ws1.Cells.Item(392, 3).Value = '14X44_ABASTO'
$ws1.Cells.Item(392, 4).Value = 59
$ws1.Cells.Item(392, 5).Value = 'LP1912'
$ws1.Cells.Item(393, 1).Value = '17:54:41'
$ws1.Cells.Item(393, 2).Value = '18:33'
$ws1.Cells.Item(393, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(393, 4).Value = 39
$ws1.Cells.Item(393, 5).Value = 'LP1912'
$ws1.Cells.Item(422, 1).Value = '18:10:23'
$ws1.Cells.Item(422, 2).Value = '19:17'
$ws1.Cells.Item(422, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(422, 4).Value = 67
$ws1.Cells.Item(422, 5).Value = 'LP1912'
$ws1.Cells.Item(423, 1).Value = '18:44:14'
$ws1.Cells.Item(423, 2).Value = '19:17'
$ws1.Cells.Item(423, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(423, 4).Value = 33
$ws1.Cells.Item(423, 5).Value = 'LP1912'
$ws1.Cells.Item(453, 1).Value = '20:11:44'
$ws1.Cells.Item(453, 2).Value = '20:11'
$ws1.Cells.Item(453, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(453, 4).Value = 0
$ws1.Cells.Item(453, 5).Value = 'LP1912'
$ws1.Cells.Item(454, 1).Value = '18:30:56'
$ws1.Cells.Item(454, 2).Value = '20:13'
$ws1.Cells.Item(454, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(454, 4).Value = 103
$ws1.Cells.Item(454, 5).Value = 'LP1912'
$ws1.Cells.Item(455, 1).Value = '20:11:44'
$ws1.Cells.Item(455, 2).Value = '20:14'
$ws1.Cells.Item(455, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(455, 4).Value = 3
$ws1.Cells.Item(455, 5).Value = 'LP1912'
$ws1.Cells.Item(456, 1).Value = '18:44:14'
$ws1.Cells.Item(456, 2).Value = '20:14'
$ws1.Cells.Item(456, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(456, 4).Value = 90
$ws1.Cells.Item(456, 5).Value = 'LP1912'
$ws1.Cells.Item(457, 1).Value = '19:47:42'
$ws1.Cells.Item(457, 2).Value = '20:15'
$ws1.Cells.Item(457, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(457, 4).Value = 28
$ws1.Cells.Item(457, 5).Value = 'LP1912'
$ws1.Cells.Item(458, 1).Value = '19:54:54'
$ws1.Cells.Item(458, 2).Value = '20:21'
$ws1.Cells.Item(458, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(458, 4).Value = 27
$ws1.Cells.Item(458, 5).Value = 'LP1912'
$ws1.Cells.Item(459, 1).Value = '20:11:44'
$ws1.Cells.Item(459, 2).Value = '20:21'
$ws1.Cells.Item(459, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(459, 4).Value = 10
$ws1.Cells.Item(459, 5).Value = 'LP1912'
$ws1.Cells.Item(460, 1).Value = '18:30:56'
$ws1.Cells.Item(460, 2).Value = '20:25'
$ws1.Cells.Item(460, 3).Value = '15_ABASTO'
$ws1.Cells.Item(460, 4).Value = 115
$ws1.Cells.Item(460, 5).Value = 'LP1912'
$ws1.Cells.Item(461, 1).Value = '18:44:14'
$ws1.Cells.Item(461, 2).Value = '20:26'
$ws1.Cells.Item(461, 3).Value = '15_ABASTO'
$ws1.Cells.Item(461, 4).Value = 102
$ws1.Cells.Item(461, 5).Value = 'LP1912'
$ws1.Cells.Item(462, 1).Value = '18:30:56'
$ws1.Cells.Item(462, 2).Value = '20:28'
$ws1.Cells.Item(462, 3).Value = '10_OLMOS'
$ws1.Cells.Item(462, 4).Value = 118
$ws1.Cells.Item(462, 5).Value = 'LP1912'
$ws1.Cells.Item(463, 1).Value = '18:44:14'
$ws1.Cells.Item(463, 2).Value = '20:29'
$ws1.Cells.Item(463, 3).Value = '10_OLMOS'
$ws1.Cells.Item(463, 4).Value = 105
$ws1.Cells.Item(463, 5).Value = 'LP1912'
$ws1.Cells.Item(464, 1).Value = '20:11:44'
$ws1.Cells.Item(464, 2).Value = '20:38'
$ws1.Cells.Item(464, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(464, 4).Value = 27
$ws1.Cells.Item(464, 5).Value = 'LP1912'
$ws1.Cells.Item(465, 1).Value = '19:11:56'
$ws1.Cells.Item(465, 2).Value = '20:43'
$ws1.Cells.Item(465, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(465, 4).Value = 92
$ws1.Cells.Item(465, 5).Value = 'LP1912'
$ws1.Cells.Item(466, 1).Value = '18:52:19'
$ws1.Cells.Item(466, 2).Value = '20:44'
$ws1.Cells.Item(466, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(466, 4).Value = 93
$ws1.Cells.Item(466, 5).Value = 'LP1912'
$ws1.Cells.Item(467, 1).Value = '18:52:19'
$ws1.Cells.Item(467, 2).Value = '20:44'
$ws1.Cells.Item(467, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(467, 4).Value = 112
$ws1.Cells.Item(467, 5).Value = 'LP1912'
$ws1.Cells.Item(468, 1).Value = '18:52:19'
$ws1.Cells.Item(468, 2).Value = '20:45'
$ws1.Cells.Item(468, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(468, 4).Value = 113
$ws1.Cells.Item(468, 5).Value = 'LP1912'
$ws1.Cells.Item(469, 1).Value = '20:11:44'
$ws1.Cells.Item(469, 2).Value = '20:46'
$ws1.Cells.Item(469, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(469, 4).Value = 35
$ws1.Cells.Item(469, 5).Value = 'LP1912'
$ws1.Cells.Item(470, 1).Value = '18:52:19'
$ws1.Cells.Item(470, 2).Value = '20:49'
$ws1.Cells.Item(470, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(470, 4).Value = 117
$ws1.Cells.Item(470, 5).Value = 'LP1912'
$ws1.Cells.Item(471, 1).Value = '19:11:56'
$ws1.Cells.Item(471, 2).Value = '20:51'
$ws1.Cells.Item(471, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(471, 4).Value = 100
$ws1.Cells.Item(471, 5).Value = 'LP1912'
$ws1.Cells.Item(472, 1).Value = '19:47:42'
$ws1.Cells.Item(472, 2).Value = '20:55'
$ws1.Cells.Item(472, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(472, 4).Value = 68
$ws1.Cells.Item(472, 5).Value = 'LP1912'
$ws1.Cells.Item(473, 1).Value = '19:11:56'
$ws1.Cells.Item(473, 2).Value = '20:56'
$ws1.Cells.Item(473, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(473, 4).Value = 105
$ws1.Cells.Item(473, 5).Value = 'LP1912'
$ws1.Cells.Item(474, 1).Value = '19:11:56'
$ws1.Cells.Item(474, 2).Value = '21:01'
$ws1.Cells.Item(474, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(474, 4).Value = 110
$ws1.Cells.Item(474, 5).Value = 'LP1912'
$ws1.Cells.Item(475, 1).Value = '19:35:19'
$ws1.Cells.Item(475, 2).Value = '21:02'
$ws1.Cells.Item(475, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(475, 4).Value = 87
$ws1.Cells.Item(475, 5).Value = 'LP1912'
$ws1.Cells.Item(476, 1).Value = '19:47:42'
$ws1.Cells.Item(476, 2).Value = '21:09'
$ws1.Cells.Item(476, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(476, 4).Value = 82
$ws1.Cells.Item(476, 5).Value = 'LP1912'
$ws1.Cells.Item(477, 1).Value = '19:47:42'
$ws1.Cells.Item(477, 2).Value = '21:23'
$ws1.Cells.Item(477, 3).Value = '10_OLMOS'
$ws1.Cells.Item(477, 4).Value = 96
$ws1.Cells.Item(477, 5).Value = 'LP1912'
$ws1.Cells.Item(478, 1).Value = '19:35:19'
$ws1.Cells.Item(478, 2).Value = '21:24'
$ws1.Cells.Item(478, 3).Value = '10_OLMOS'
$ws1.Cells.Item(478, 4).Value = 109
$ws1.Cells.Item(478, 5).Value = 'LP1912'
$ws1.Cells.Item(479, 1).Value = '19:54:54'
$ws1.Cells.Item(479, 2).Value = '21:48'
$ws1.Cells.Item(479, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(479, 4).Value = 114
$ws1.Cells.Item(479, 5).Value = 'LP1912'
$ws1.Cells.Item(480, 1).Value = '20:11:44'
$ws1.Cells.Item(480, 2).Value = '21:49'
$ws1.Cells.Item(480, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(480, 4).Value = 98
$ws1.Cells.Item(480, 5).Value = 'LP1912'
$ws1.Cells.Item(481, 1).Value = '20:11:44'
$ws1.Cells.Item(481, 2).Value = '21:55'
$ws1.Cells.Item(481, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(481, 4).Value = 104
$ws1.Cells.Item(481, 5).Value = 'LP1912'

# Sheet 2 and 3 timestamp refresh
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 20:11:44'
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 20:11:44'
